$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.013.98'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '2.421.32'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.74'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.18'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('E9').Value = '  -7.47%  '
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('E12').Value = '  -4.01%  '
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.20'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('D16').Value = '2.841.66'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = '61.996.54'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').Value = '2.431.92'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.30'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '323.61'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('E21').Value = '  +0.87%  '
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.07'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.56%  '
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.76'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '555.28'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -5.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('E30').Value = '  -1.41%  '
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('E32').Value = '  -4.48%  '
$ws.Range('E33').Value = '  -1.87%  '
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.51'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.46%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.75'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('E38').Value = '  -1.29%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '152.99'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.45'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.61%  '
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('E42').Value = '  -1.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.995'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '147.38'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.23'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.54%  '
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0527'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.05%  '
$ws.Range('E48').Value = '  +0.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.81'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0918'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('E51').Value = '  -0.42%  '

Write-Output "Updated cryptos list"
